$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 0.03675
    "B5"  = 5.524925000000001
    "B9"  = 0.1716
    "B13" = 2.2212
    "B15" = 12.683135
    "B16" = 1.904905
    "B18" = 2.822072
    "B19" = 0.126225
    "B20" = 11.7016
    "B23" = 2.34657
    "B24" = 3.46575
    "B27" = 0.006000000000000002
    "B29" = 2.021495
    "B30" = 0.425425
    "B33" = 0.8747999999999999
    "B39" = 0.277725
    "B44" = 6.219109
    "B45" = 0.37315
    "B48" = 1.86171
    "B52" = 4.041088
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
